$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.484.46"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "'1.910.99"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'327.56"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4745"
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("D8").Value = "'0.4094"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'47.70"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "'0.08031"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'1.009"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "'22.41"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "'1.911.30"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'5.951"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'7.155"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "'89.32"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.00001031"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.06596"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'17.74"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'29.503.65"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "'5.549"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").Value = "'2.143.56"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'153.56"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").Value = "'19.77"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'5.750"
$ws.Range("E29").Value = "  +6.48%  "
$ws.Range("D30").Value = "'2.135"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "'117.50"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'1.068"
$ws.Range("E32").Value = "  +9.31%  "
$ws.Range("D33").Value = "'0.09556"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "'1.422"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "'5.387"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06089"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02253"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "'8.363"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'1.174"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "'0.5882"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "'0.1842"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'10.14"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "'1.302"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "'2.409"
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("D46").Value = "'0.07786"
$ws.Range("E46").Value = "  +10.89%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.23"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5543"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'1.933"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "'113.43"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'44.17"
$ws.Range("E51").Value = "  -8.14%  "
